$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Inscritos (E), new Pagos (F), new Inscricoes homologadas (H)
# (Isencoes deferidas / column G is unchanged at 0 throughout)
$updates = @(
    @(2, 71, 42, 42),
    @(3, 25, 20, 20),
    @(5, 81, 45, 45),
    @(6, 27, 10, 10),
    @(7, 22, 10, 10),
    @(9, 7, 5, 5),
    @(10, 244, 107, 107),
    @(11, 179, 94, 94),
    @(12, 267, 132, 132),
    @(13, 84, 39, 39),
    @(14, 75, 35, 35),
    @(16, 109, 54, 54),
    @(17, 51, 23, 23),
    @(18, 39, 16, 16),
    @(20, 61, 21, 21),
    @(21, 82, 44, 44),
    @(22, 105, 49, 49),
    @(23, 112, 52, 52),
    @(24, 121, 55, 55),
    @(25, 120, 51, 51),
    @(26, 76, 41, 41),
    @(27, 174, 87, 87),
    @(28, 105, 32, 32),
    @(29, 112, 61, 61),
    @(30, 126, 65, 65),
    @(31, 45, 21, 21),
    @(32, 117, 59, 59),
    @(33, 163, 70, 70),
    @(34, 124, 69, 69),
    @(35, 89, 47, 47),
    @(37, 88, 40, 40),
    @(39, 120, 49, 49),
    @(40, 157, 64, 64),
    @(41, 215, 85, 85),
    @(42, 201, 94, 94),
    @(43, 63, 27, 27),
    @(44, 169, 79, 79),
    @(45, 69, 34, 34),
    @(46, 152, 76, 76),
    @(47, 250, 110, 110),
    @(48, 120, 41, 41),
    @(49, 140, 57, 57),
    @(50, 113, 45, 45),
    @(51, 116, 46, 46),
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 5).Value = $u[1]
    $ws.Cells.Item($row, 6).Value = $u[2]
    $ws.Cells.Item($row, 8).Value = $u[3]
}
